# The commit swaps the East-Asian fallback font from "DejaVu Sans" to
# "Tahoma" for the document's base styles, and makes the complex-script
# ("cs") font explicit ("DejaVu Sans") on a few styles that previously
# inherited it implicitly (empty <w:rPr/>).
#
# rFonts <-> Font property mapping used below:
#   w:ascii            -> Font.NameAscii (Font.Name sets ascii+hAnsi both)
#   w:eastAsia         -> Font.NameFarEast
#   w:cs (complex script) -> Font.NameBi  (NameBidirectional)

$d = $word.ActiveDocument

# "Normal" paragraph style: eastAsia DejaVu Sans -> Tahoma
$normal = $d.Styles.Item("Normal")
$normal.Font.NameFarEast = "Tahoma"

# "Heading" paragraph style: eastAsia DejaVu Sans -> Tahoma
$heading = $d.Styles.Item("Heading")
$heading.Font.NameFarEast = "Tahoma"

# "List" paragraph style: make the complex-script font explicit
$list = $d.Styles.Item("List")
$list.Font.NameBi = "DejaVu Sans"

# "Caption" paragraph style: make the complex-script font explicit
$caption = $d.Styles.Item("Caption")
$caption.Font.NameBi = "DejaVu Sans"

# "Index" paragraph style: make the complex-script font explicit
$index = $d.Styles.Item("Index")
$index.Font.NameBi = "DejaVu Sans"
